$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4 to become the "charging_station" entry
$ws.Range("A4").Value = "charging_station1"
$ws.Range("B4").Value = "charging_station"

# Remove the old row 5 (previously "bat1"/"bat") entirely
$ws.Rows("5:5").Delete()
